$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1333.4546
$ws.Range("I38").Value = 1066.8572
$ws.Range("J38").Value = 1800
$ws.Range("K38").Value = 3200.5716
$ws.Range("L38").Value = 5400
$ws.Range("M38").Value = -2828.5716
$ws.Range("N38").Value = -6144
$ws.Range("H82").Value = 3447220.5
$ws.Range("I82").Value = 4020090.8
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 12060272.4
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -12059866.4
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 3447220.5
$ws.Range("I85").Value = 4020090.8
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 12060272.4
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -12058868.4
$ws.Range("N85").Value = -32808
$ws.Range("H107").Value = 564.86664
$ws.Range("I107").Value = 564.7778
$ws.Range("J107").Value = 565
$ws.Range("K107").Value = 564.7778
$ws.Range("L107").Value = 565
$ws.Range("M107").Value = 1355.2222
$ws.Range("N107").Value = -4405
$ws.Range("H112").Value = 1394.8572
$ws.Range("I112").Value = 633.3333
$ws.Range("J112").Value = 1521.7778
$ws.Range("K112").Value = 1899.9999
$ws.Range("L112").Value = 4565.3334
$ws.Range("M112").Value = -791.9999
$ws.Range("N112").Value = -6781.3334
$ws.Range("H113").Value = 1525.2106
$ws.Range("I113").Value = 1562.25
$ws.Range("J113").Value = 1515.3334
$ws.Range("K113").Value = 1562.25
$ws.Range("L113").Value = 1515.3334
$ws.Range("M113").Value = 1691.75
$ws.Range("N113").Value = -8023.3334
$ws.Range("H116").Value = 3500407.8
$ws.Range("I116").Value = 25643590
$ws.Range("J116").Value = 4115.7896
$ws.Range("K116").Value = 25643590
$ws.Range("L116").Value = 4115.7896
$ws.Range("M116").Value = -25640148
$ws.Range("N116").Value = -10999.7896
$ws.Range("H125").Value = 3959.4285
$ws.Range("I125").Value = 1980
$ws.Range("J125").Value = 4289.3335
$ws.Range("K125").Value = 17820
$ws.Range("L125").Value = 38604.0015
$ws.Range("M125").Value = -15360
$ws.Range("N125").Value = -43524.0015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 33900
$ws.Range("J52").Value = 33900
$ws.Range("L52").Value = 33900
$ws.Range("N52").Value = -34536
$ws.Range("H61").Value = 2003.238
$ws.Range("I61").Value = 1296.0667
$ws.Range("J61").Value = 3771.1667
$ws.Range("K61").Value = 1296.0667
$ws.Range("L61").Value = 3771.1667
$ws.Range("M61").Value = -1084.0667
$ws.Range("N61").Value = -4195.1667
$ws.Range("H74").Value = 1406.5278
$ws.Range("I74").Value = 1386.4642
$ws.Range("J74").Value = 1476.75
$ws.Range("K74").Value = 1386.4642
$ws.Range("L74").Value = 1476.75
$ws.Range("M74").Value = -512.4641999999999
$ws.Range("N74").Value = -3224.75
$ws.Range("H77").Value = 1406.5278
$ws.Range("I77").Value = 1386.4642
$ws.Range("J77").Value = 1476.75
$ws.Range("K77").Value = 6932.321
$ws.Range("L77").Value = 7383.75
$ws.Range("M77").Value = -2564.321
$ws.Range("N77").Value = -16119.75
$ws.Range("H119").Value = 22474
$ws.Range("J119").Value = 22474
$ws.Range("L119").Value = 22474
$ws.Range("N119").Value = -32150
$ws.Range("H122").Value = 1471.7084
$ws.Range("I122").Value = 1402.7333
$ws.Range("J122").Value = 1586.6666
$ws.Range("K122").Value = 4208.199900000001
$ws.Range("L122").Value = 4759.9998
$ws.Range("M122").Value = -1758.199900000001
$ws.Range("N122").Value = -9659.9998
$ws.Range("H136").Value = 2003.238
$ws.Range("I136").Value = 1296.0667
$ws.Range("J136").Value = 3771.1667
$ws.Range("K136").Value = 3888.2001
$ws.Range("L136").Value = 11313.5001
$ws.Range("M136").Value = -1338.2001
$ws.Range("N136").Value = -16413.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1975.5
$ws.Range("J22").Value = 1975.5
$ws.Range("L22").Value = 5926.5
$ws.Range("N22").Value = -6264.5
$ws.Range("H27").Value = 1975.5
$ws.Range("J27").Value = 1975.5
$ws.Range("L27").Value = 5926.5
$ws.Range("N27").Value = -6130.5
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H131").Value = 758.61365
$ws.Range("I131").Value = 442.30768
$ws.Range("J131").Value = 891.25806
$ws.Range("K131").Value = 1326.92304
$ws.Range("L131").Value = 2673.77418
$ws.Range("M131").Value = 3713.07696
$ws.Range("N131").Value = -12753.77418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3333.3333
$ws.Range("I102").Value = 3333.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3333.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1711.3333
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 14287529
$ws.Range("I122").Value = 33334566
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 100003698
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -100001248
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 559
$ws.Range("I39").Value = 559
$ws.Range("K39").Value = 559
$ws.Range("M39").Value = -99
$ws.Range("H69").Value = 10019900
$ws.Range("J69").Value = 10019900
$ws.Range("L69").Value = 10019900
$ws.Range("N69").Value = -10021522
$ws.Range("H72").Value = 10019900
$ws.Range("J72").Value = 10019900
$ws.Range("L72").Value = 30059700
$ws.Range("N72").Value = -30067812
$ws.Range("H127").Value = 62500
$ws.Range("J127").Value = 62500
$ws.Range("L127").Value = 62500
$ws.Range("N127").Value = -72420
$ws.Range("H136").Value = 2563.8262
$ws.Range("I136").Value = 2185.1875
$ws.Range("J136").Value = 3429.2856
$ws.Range("K136").Value = 6555.5625
$ws.Range("L136").Value = 10287.8568
$ws.Range("M136").Value = -4005.5625
$ws.Range("N136").Value = -15387.8568
